# Edit described by commit: "added upper case option and multiple receivers supported"
#
# The single slide has a banner ("LA MULTI ANI") and a lower textbox that used
# to be an empty placeholder (formatted for ALL-CAPS display via cap="all").
# The edit fills that textbox with a list of honorees (the "receivers"),
# typed directly in upper case (the "upper case option"), one per paragraph
# (the "multiple receivers" support), and resizes/repositions the box to fit
# the new, larger block of text.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shp = $s.Shapes.Item(2)

# Shape was "TextBox 3" -- rename + reflow / resize to its new footprint.
$shp.Name = "TextBox 2"
$shp.Left = 56.38551181102362
$shp.Top = 129.13236220472442
$shp.Width = 845.7831496062992
$shp.Height = 383.338188976378

$tf = $shp.TextFrame
$tr = $tf.TextRange

$names = @(
    "SUCIU IOAN",
    "ROSCA MARIA",
    "FINDEIS CRISTINA",
    "BRANDZANIC MARCEL - IVAN",
    "MUSKA MIHAELA",
    "JENTIMIR LAURA",
    "POPA NARCIS ADRIAN",
    "TATARU BEATRICE"
)

# First paragraph
$tr.Text = $names[0]
$tr.Font.Size = 80
$tr.Font.Bold = $true
$tr.LanguageID = "en-US"

# Remaining paragraphs - append as new paragraphs and (re)apply formatting to
# the whole range each time so every run picks up size/bold/language.
for ($i = 1; $i -lt $names.Count; $i++) {
    [void]$tr.InsertAfter("`r" + $names[$i])
    $tr2 = $tf.TextRange
    $tr2.Font.Size = 80
    $tr2.Font.Bold = $true
    $tr2.LanguageID = "en-US"
}

Write-Host "Updated" $tf.TextRange.Paragraphs().Count "paragraphs in" $shp.Name
